{"js": "// Apply three text edits to the body of the document using the Word\n// JavaScript API (Office.js). Each edit is located with a narrow, unique\n// search string so that Body.search() matches exactly one range.\n\nconst body = context.document.body;\n\n// ---------------------------------------------------------------------\n// Edit 1: Insert a new sentence into the \"However, the main difference...\"\n// paragraph, right after \"...prove fairness in clustering algorithms.\"\n// and before \" This is likely due to its clean nature...\".\n// ---------------------------------------------------------------------\n{\n  const results = body.search(\n    \"prove fairness in clustering algorithms. This is likely due to\",\n    { matchCase: true, matchWholeWord: false }\n  );\n  results.load(\"items/text\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    results.items[0].insertText(\n      \"prove fairness in clustering algorithms. Thus, the state-of-the-art method was more akin to clustering over classification. This is likely due to\",\n      \"Replace\"\n    );\n  }\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// Edit 2: Remove the stray \" m\" typo: \"algorithm m in handling\" ->\n// \"algorithm in handling\".\n// ---------------------------------------------------------------------\n{\n  const results = body.search(\"effectiveness of their algorithm m in handling\", {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items/text\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    results.items[0].insertText(\n      \"effectiveness of their algorithm in handling\",\n      \"Replace\"\n    );\n  }\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// Edit 3: Replace the entire \"The heading of a section should be...\"\n// paragraph body text with a single period.\n// ---------------------------------------------------------------------\n{\n  const results = body.search(\n    \"The heading of a section should be in Times New Roman 12-point bold in all-capitals flush left with an additional 6-points of white space above the section head.  Sections and subsequent sub- sections should be numbered and flush left. For a section head and a subsection head together (such as Section 3 and subsection 3.1), use no additional space above the subsection head.\",\n    { matchCase: true, matchWholeWord: false }\n  );\n  results.load(\"items/text\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    results.items[0].insertText(\".\", \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply three text edits to the document body using the Word COM object\n# model (Find/Replace). Each Find.Text value is a long, unique substring\n# so only the intended occurrence is touched, and wdReplaceOne (2) keeps\n# each call to a single replacement.\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# Edit 1: Insert a new sentence into the \"However, the main difference...\"\n# paragraph, right after \"...prove fairness in clustering algorithms.\"\n# and before \" This is likely due to its clean nature...\".\n# ---------------------------------------------------------------------\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Text = \"prove fairness in clustering algorithms. This is likely due to\"\n$find1.Replacement.Text = \"prove fairness in clustering algorithms. Thus, the state-of-the-art method was more akin to clustering over classification. This is likely due to\"\n$find1.Execute(\n    $find1.Text,\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    $find1.Replacement.Text,\n    2\n)\n\n# ---------------------------------------------------------------------\n# Edit 2: Remove the stray \" m\" typo: \"algorithm m in handling\" ->\n# \"algorithm in handling\".\n# ---------------------------------------------------------------------\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"effectiveness of their algorithm m in handling\"\n$find2.Replacement.Text = \"effectiveness of their algorithm in handling\"\n$find2.Execute(\n    $find2.Text,\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    $find2.Replacement.Text,\n    2\n)\n\n# ---------------------------------------------------------------------\n# Edit 3: Replace the entire \"The heading of a section should be...\"\n# paragraph body text with a single period.\n# ---------------------------------------------------------------------\n$find3 = $d.Content.Find\n$find3.ClearFormatting()\n$find3.Replacement.ClearFormatting()\n$find3.Text = \"The heading of a section should be in Times New Roman 12-point bold in all-capitals flush left with an additional 6-points of white space above the section head.  Sections and subsequent sub- sections should be numbered and flush left. For a section head and a subsection head together (such as Section 3 and subsection 3.1), use no additional space above the subsection head.\"\n$find3.Replacement.Text = \".\"\n$find3.Execute(\n    $find3.Text,\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    $find3.Replacement.Text,\n    2\n)\n"}
